$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D cells to remain text even when the new value looks numeric,
# matching the original inline-string (text) storage used by the source file.

$ws.Range("D2").Value = '71.869.75'
$ws.Range("E2").Value = '  +3.47%  '
$ws.Range("D3").Value = '3.634.15'
$ws.Range("E3").Value = '  +7.19%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '588.33'
$ws.Range("E5").Value = '  +0.65%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '181.14'
$ws.Range("E6").Value = '  +0.54%  '
$ws.Range("D7").Value = '3.627.50'
$ws.Range("E7").Value = '  +7.23%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.616'
$ws.Range("E8").Value = '  +3.24%  '
$ws.Range("E9").Value = '  +0.10%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.203'
$ws.Range("E10").Value = '  +1.30%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.608'
$ws.Range("E11").Value = '  +2.94%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '49.68'
$ws.Range("E12").Value = '  +2.98%  '
$ws.Range("E13").Value = '  +0.37%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '684.44'
$ws.Range("E14").Value = '  +0.28%  '
$ws.Range("D15").Value = '4.220.34'
$ws.Range("E15").Value = '  +6.72%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '9.03'
$ws.Range("E16").Value = '  +4.44%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '71.835.77'
$ws.Range("E17").Value = '  +3.24%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.633.04'
$ws.Range("E18").Value = '  +6.53%  '
$ws.Range("E19").Value = '  +1.81%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.32'
$ws.Range("E20").Value = '  +3.50%  '
$ws.Range("E21").Value = '  +2.63%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.939'
$ws.Range("E22").Value = '  +3.31%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.91'
$ws.Range("E23").Value = '  +10.48%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '17.81'
$ws.Range("E24").Value = '  +3.22%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '103.36'
$ws.Range("E25").Value = '  +0.85%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.02'
$ws.Range("E26").Value = '  +2.36%  '
$ws.Range("E27").Value = '  +5.21%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.99'
$ws.Range("E28").Value = '  +2.94%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '35.22'
$ws.Range("E29").Value = '  +4.04%  '
$ws.Range("E30").Value = '  +4.87%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.36'
$ws.Range("E31").Value = '  +6.04%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.20'
$ws.Range("E32").Value = '  +16.48%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '586.26'
$ws.Range("E33").Value = '  +5.83%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '11.31'
$ws.Range("E34").Value = '  +1.95%  '
$ws.Range("E35").Value = '  +2.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '59.36'
$ws.Range("E36").Value = '  +1.72%  '
$ws.Range("E37").Value = '  +0.00%  '
$ws.Range("D38").Value = '3.674.98'
$ws.Range("E38").Value = '  +0.37%  '
$ws.Range("E39").Value = '  +1.86%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '35.69'
$ws.Range("E40").Value = '  +0.13%  '
$ws.Range("D41").Value = '0.0₃0764'
$ws.Range("E41").Value = '  +5.05%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.42'
$ws.Range("E42").Value = '  +3.83%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0468'
$ws.Range("E43").Value = '  +9.17%  '
$ws.Range("E44").Value = '  +2.72%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.345'
$ws.Range("E45").Value = '  +2.29%  '
$ws.Range("E46").Value = '  +1.63%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.81'
$ws.Range("E47").Value = '  +5.32%  '
$ws.Range("E48").Value = '  +3.28%  '
$ws.Range("E49").Value = '  +3.61%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.998'
$ws.Range("E50").Value = '  -0.35%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '131.59'
$ws.Range("E51").Value = '  +1.81%  '
